$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a number (e.g. '309.68') need the
# column forced to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value and the grouped-dot price strings
# (e.g. '44.631.79') would lose their literal text representation.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '44.631.79'
$ws.Range('D3').Value = '2.251.95'
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').Value = '309.68'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').Value = '95.78'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').Value = '35.35'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  +0.99%  '
$ws.Range('D12').Value = '7.34'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.845'
$ws.Range('E14').Value = '  +3.70%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.248.67'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '13.70'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '44.260.40'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '0.0₃0968'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').Value = '6.42'
$ws.Range('E19').Value = '  +4.30%  '
$ws.Range('D20').Value = '12.24'
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '65.92'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D22').Value = '240.09'
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('E23').Value = '  +3.24%  '
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  +4.29%  '
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  +5.25%  '
$ws.Range('D27').Value = '9.90'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').Value = '37.74'
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('D29').Value = '6.07'
$ws.Range('E29').Value = '  +3.57%  '
$ws.Range('D30').Value = '20.20'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = '0.0810'
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('D34').Value = '3.17'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('D36').Value = '0.121'
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('D37').Value = '1.83'
$ws.Range('E37').Value = '  +3.59%  '
$ws.Range('D38').Value = '3.47'
$ws.Range('E38').Value = '  +4.59%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '3.83'
$ws.Range('E39').Value = '  +1.61%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '14.47'
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('D43').Value = '1.754.91'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Value = '0.195'
$ws.Range('E44').Value = '  +5.49%  '
$ws.Range('D45').Value = '81.45'
$ws.Range('E45').Value = '  -3.87%  '
$ws.Range('D46').Value = '71.39'
$ws.Range('E46').Value = '  +4.52%  '
$ws.Range('D47').Value = '100.11'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').Value = '55.83'
$ws.Range('E48').Value = '  +3.49%  '
$ws.Range('D49').Value = '8.22'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').Value = '4.89'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('E51').Value = '  +5.14%  '
